$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Tests" (sheet1.xml)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Tests")

# New column I ("param:scope") with width 20 (XML width=20 <=> ColumnWidth=19.17)
$ws1.Columns.Item(9).ColumnWidth = 19.17

# Header for new column
$ws1.Range("I1").Value = "param:scope"

# New value for existing row 2
$ws1.Range("I2").Value = "work"

# New row 3 - "get-facets - Missing Required Param" test case
$ws1.Range("A3").Value = "get-facets - Missing Required Param"
$ws1.Range("B3").Value = "Test GET /api/facets/:scope with missing required parameters"
# leading apostrophe forces literal text "true" instead of a boolean TRUE
$ws1.Range("C3").Value = "'true"
$ws1.Range("D3").Value = 400
$ws1.Range("E3").Value = 10000
$ws1.Range("F3").Value = 2000
$ws1.Range("G3").Value = 500
$ws1.Range("H3").Value = "get-facets,validation"
# empty string cell (apostrophe forces an explicit empty text value rather than a blank cell)
$ws1.Range("I3").Value = "'"

# ---------------------------------------------------------------------------
# Sheet "Documentation" (sheet2.xml)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Documentation")

# Insert two new rows at position 18, pushing the "Endpoint-Specific Notes"
# block (and everything below it) down by two rows.
$ws2.Rows.Item(18).Insert()
$ws2.Rows.Item(18).Insert()

# Row 17 gains a heading
$ws2.Range("A17").Value = "Parameter Descriptions:"

# New rows 18-19 (parameter description + blank separator)
$ws2.Range("A18").Value = "param:scope"
$ws2.Range("B18").Value = "Search scope (work, person, place, concept, event, etc.) (string) (REQUIRED - highlighted in yellow)"
$ws2.Range("A19").Value = "'"

# New row 24 at the end, documenting required parameters
$ws2.Range("A24").Value = "• Required parameters: scope"
